# Included PCE results in the log outputs and details on dataset size.
#
# The "feature importance" table (A2:B30) is re-ranked / refreshed with the
# new importances that resulted from incorporating the PCE results into the
# training run. Row 31 (COSMO_Screening_Charge / 0) and the header row are
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "deltaE_RedOxH"
$ws.Range("B2").Value = 0.0852532962900039
$ws.Range("A3").Value = "IP"
$ws.Range("B3").Value = 0.06734353704942037
$ws.Range("A4").Value = "HOMO"
$ws.Range("B4").Value = 0.05530325205805955
$ws.Range("A5").Value = "RotatableBonds"
$ws.Range("B5").Value = 0.05459798518767987
$ws.Range("A6").Value = "deltaE_LCB"
$ws.Range("B6").Value = 0.04825404361252158
$ws.Range("A7").Value = "electronegativity"
$ws.Range("B7").Value = 0.04680774305276485
$ws.Range("A8").Value = "Max_f_osc"
$ws.Range("B8").Value = 0.04525695612708441
$ws.Range("A9").Value = "electrophilicityIndex"
$ws.Range("B9").Value = 0.0449889657122557
$ws.Range("A10").Value = "Surface_Area_A2"
$ws.Range("B10").Value = 0.04028816928167109
$ws.Range("A11").Value = "Molecular_Volume_A3"
$ws.Range("B11").Value = 0.03867154047465479
$ws.Range("A12").Value = "HBondAcceptors"
$ws.Range("B12").Value = 0.03714603299851457
$ws.Range("A13").Value = "deltaE_HL"
$ws.Range("B13").Value = 0.03569813232664815
$ws.Range("A14").Value = "LHE"
$ws.Range("B14").Value = 0.03287102828368097
$ws.Range("A15").Value = "Total_Energy_Hartree"
$ws.Range("B15").Value = 0.03196260193340262
$ws.Range("A16").Value = "electroacceptingPower"
$ws.Range("B16").Value = 0.03097894467992434
$ws.Range("A17").Value = "TPSA"
$ws.Range("B17").Value = 0.03000586522158782
$ws.Range("A18").Value = "chemHardness"
$ws.Range("B18").Value = 0.0299488781731147
$ws.Range("A19").Value = "RingCount"
$ws.Range("B19").Value = 0.02631034166682851
$ws.Range("A20").Value = "electrodonatingPower"
$ws.Range("B20").Value = 0.02476738987344099
$ws.Range("A21").Value = "HBondDonors"
$ws.Range("B21").Value = 0.02443145326078758
$ws.Range("A22").Value = "LUMO"
$ws.Range("B22").Value = 0.02344892064026551
$ws.Range("A23").Value = "Max_Absorption_nm"
$ws.Range("B23").Value = 0.02273387369672214
$ws.Range("A24").Value = "EA"
$ws.Range("B24").Value = 0.02182879962396833
$ws.Range("A25").Value = "AromaticRings"
$ws.Range("B25").Value = 0.02050164763705128
$ws.Range("A26").Value = "LogP"
$ws.Range("B26").Value = 0.01946073814220173
$ws.Range("A27").Value = "elnChemPot"
$ws.Range("B27").Value = 0.01839312867180041
$ws.Range("A28").Value = "Dipole_Moment"
$ws.Range("B28").Value = 0.01493344666395918
$ws.Range("A29").Value = "Solvation_Energy_eV"
$ws.Range("B29").Value = 0.01399992049796566
$ws.Range("A30").Value = "Mass"
$ws.Range("B30").Value = 0.01381336716201939
